$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (source row 85)
$ws.Range("D2").Value = 44400
$ws.Range("J2").Value = 240
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 15000
$ws.Range("P2").Value = 1500

# Row 3 (source row 146)
$ws.Range("D3").Value = 44239
$ws.Range("J3").Value = 250

# Row 4 (source row 141)
$ws.Range("D4").Value = 44162
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 12000
$ws.Range("P4").Value = 1200

# Row 5 (source row 12)
$ws.Range("D5").Value = 44305
$ws.Range("J5").Value = 40

# Row 6 (source row 38)
$ws.Range("D6").Value = 44348
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("P6").Value = 1500

# Row 7 (source row 93)
$ws.Range("D7").Value = 44285
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 14000
$ws.Range("M7").Value = 14400
$ws.Range("P7").Value = 1440

# Row 8 (source row 29)
$ws.Range("D8").Value = 44329
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("P8").Value = 1500

# Row 9 (source row 19)
$ws.Range("D9").Value = 44189
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 13500
$ws.Range("M9").Value = 13250
$ws.Range("P9").Value = 1325

# Row 10 (source row 111)
$ws.Range("D10").Value = 44267
$ws.Range("J10").Value = 250
$ws.Range("M10").Value = 14600
$ws.Range("P10").Value = 1460

# Row 11 (source row 28)
$ws.Range("D11").Value = 44396
$ws.Range("J11").Value = 40
$ws.Range("K11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("P11").Value = 1500

# Row 12 (source row 41)
$ws.Range("D12").Value = 44313
$ws.Range("J12").Value = 240

# Row 13 (source row 2)
$ws.Range("D13").Value = 44181
$ws.Range("J13").Value = 30
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 12000
$ws.Range("P13").Value = 1200

# Row 14 (source row 107)
$ws.Range("D14").Value = 44249
$ws.Range("J14").Value = 60

# Row 15 (source row 98)
$ws.Range("D15").Value = 44266
$ws.Range("K15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 15000
$ws.Range("P15").Value = 1500

# Row 16 (source row 145)
$ws.Range("D16").Value = 44176
$ws.Range("J16").Value = 250
$ws.Range("K16").Value = 12000
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = 12000
$ws.Range("P16").Value = 1200

# Row 17 (source row 75)
$ws.Range("D17").Value = 44383
$ws.Range("J17").Value = 240
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 15000
$ws.Range("P17").Value = 1500

# Row 18 (source row 90)
$ws.Range("D18").Value = 44208
$ws.Range("J18").Value = 250
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 14000
$ws.Range("P18").Value = 1400

# Row 19 (source row 129)
$ws.Range("D19").Value = 44277
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 14000
$ws.Range("P19").Value = 1400

# Row 20 (source row 87)
$ws.Range("D20").Value = 44442
$ws.Range("J20").Value = 240
$ws.Range("K20").Value = 18000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 18000
$ws.Range("P20").Value = 1800

# Row 21 (source row 92)
$ws.Range("D21").Value = 44354
$ws.Range("J21").Value = 40

# Row 22 (source row 27)
$ws.Range("D22").Value = 44179
$ws.Range("K22").Value = 12000
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = 12000
$ws.Range("P22").Value = 1200

# Row 23 (source row 34)
$ws.Range("D23").Value = 44302
$ws.Range("K23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 15000
$ws.Range("P23").Value = 1500

# Row 24 (source row 88)
$ws.Range("D24").Value = 44309
$ws.Range("J24").Value = 250

# Row 25 (source row 62)
$ws.Range("D25").Value = 44253
$ws.Range("J25").Value = 250
$ws.Range("K25").Value = 15000
$ws.Range("M25").Value = 15000
$ws.Range("P25").Value = 1500

# Row 26 (source row 152)
$ws.Range("D26").Value = 44323
$ws.Range("J26").Value = 250

# Row 27 (source row 22)
$ws.Range("D27").Value = 44462
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 17000
$ws.Range("M27").Value = 17000
$ws.Range("P27").Value = 1700

# Row 28 (source row 67)
$ws.Range("D28").Value = 44342
$ws.Range("J28").Value = 25

# Row 29 (source row 13)
$ws.Range("D29").Value = 44194
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 13000
$ws.Range("L29").Value = 13500
$ws.Range("M29").Value = 13250
$ws.Range("P29").Value = 1325

# Row 30 (source row 83)
$ws.Range("D30").Value = 44341
$ws.Range("J30").Value = 250

# Row 31 (source row 133)
$ws.Range("D31").Value = 44160
$ws.Range("J31").Value = 27
$ws.Range("K31").Value = 12000
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = 12000
$ws.Range("P31").Value = 1200

# Row 32 (source row 118)
$ws.Range("D32").Value = 44196
$ws.Range("J32").Value = 300
$ws.Range("K32").Value = 13000
$ws.Range("L32").Value = 13500
$ws.Range("M32").Value = 13250
$ws.Range("P32").Value = 1325

# Row 33 (source row 101)
$ws.Range("D33").Value = 44460
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 17000
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = 17000
$ws.Range("P33").Value = 1700

# Row 34 (source row 125)
$ws.Range("D34").Value = 44371
$ws.Range("J34").Value = 80

# Row 35 (source row 55)
$ws.Range("D35").Value = 44299
$ws.Range("J35").Value = 250

# Row 36 (source row 94)
$ws.Range("D36").Value = 44167
$ws.Range("J36").Value = 20
$ws.Range("K36").Value = 12000
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 12000
$ws.Range("P36").Value = 1200

# Row 37 (source row 97)
$ws.Range("D37").Value = 44217
$ws.Range("J37").Value = 150
$ws.Range("M37").Value = 14533
$ws.Range("P37").Value = 1453

# Row 38 (source row 135)
$ws.Range("D38").Value = 44365

# Row 39 (source row 60)
$ws.Range("D39").Value = 44340
$ws.Range("J39").Value = 40
$ws.Range("K39").Value = 15000
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = 15000
$ws.Range("P39").Value = 1500

# Row 40 (source row 52)
$ws.Range("D40").Value = 44238
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 15000
$ws.Range("M40").Value = 15000
$ws.Range("P40").Value = 1500

# Row 41 (source row 147)
$ws.Range("D41").Value = 44376

# Row 42 (source row 53)
$ws.Range("D42").Value = 44446
$ws.Range("J42").Value = 240
$ws.Range("K42").Value = 18000
$ws.Range("L42").Value = 18000
$ws.Range("M42").Value = 18000
$ws.Range("P42").Value = 1800

# Row 43 (source row 25)
$ws.Range("D43").Value = 44222
$ws.Range("J43").Value = 300
$ws.Range("K43").Value = 14000
$ws.Range("M43").Value = 14500
$ws.Range("P43").Value = 1450

# Row 44 (source row 24)
$ws.Range("D44").Value = 44379
$ws.Range("J44").Value = 240
$ws.Range("K44").Value = 15000
$ws.Range("M44").Value = 15000
$ws.Range("P44").Value = 1500

# Row 45 (source row 110)
$ws.Range("D45").Value = 44260
$ws.Range("J45").Value = 300
$ws.Range("K45").Value = 14500
$ws.Range("L45").Value = 15000
$ws.Range("M45").Value = 14750
$ws.Range("P45").Value = 1475

# Row 46 (source row 150)
$ws.Range("D46").Value = 44211
$ws.Range("J46").Value = 240
$ws.Range("K46").Value = 14000
$ws.Range("M46").Value = 14500
$ws.Range("P46").Value = 1450

# Row 47 (source row 74)
$ws.Range("D47").Value = 44218
$ws.Range("J47").Value = 250
$ws.Range("K47").Value = 14000
$ws.Range("L47").Value = 14000
$ws.Range("M47").Value = 14000
$ws.Range("P47").Value = 1400

# Row 48 (source row 95)
$ws.Range("D48").Value = 44210
$ws.Range("K48").Value = 14000
$ws.Range("M48").Value = 14500
$ws.Range("P48").Value = 1450

# Row 49 (source row 144)
$ws.Range("D49").Value = 44257
$ws.Range("J49").Value = 250

# Row 50 (source row 89)
$ws.Range("D50").Value = 44200
$ws.Range("J50").Value = 80
$ws.Range("K50").Value = 13000
$ws.Range("L50").Value = 13000
$ws.Range("M50").Value = 13000
$ws.Range("P50").Value = 1300

# Row 51 (source row 59)
$ws.Range("D51").Value = 44417
$ws.Range("J51").Value = 80

# Row 52 (source row 77)
$ws.Range("D52").Value = 44246
$ws.Range("J52").Value = 250

# Row 53 (source row 4)
$ws.Range("D53").Value = 44295
$ws.Range("J53").Value = 250
$ws.Range("K53").Value = 13000
$ws.Range("L53").Value = 14000
$ws.Range("M53").Value = 13400
$ws.Range("P53").Value = 1340

# Row 54 (source row 68)
$ws.Range("D54").Value = 44225
$ws.Range("H54").Value = "Chino"
$ws.Range("J54").Value = 300
$ws.Range("K54").Value = 14000
$ws.Range("L54").Value = 15000
$ws.Range("M54").Value = 14500
$ws.Range("P54").Value = 1450

# Row 55 (source row 72)
$ws.Range("D55").Value = 44399
$ws.Range("J55").Value = 120

# Row 56 (source row 70)
$ws.Range("D56").Value = 44327
$ws.Range("J56").Value = 250

# Row 57 (source row 138)
$ws.Range("D57").Value = 44175
$ws.Range("J57").Value = 150
$ws.Range("K57").Value = 12000
$ws.Range("L57").Value = 12000
$ws.Range("M57").Value = 12000
$ws.Range("P57").Value = 1200

# Row 58 (source row 30)
$ws.Range("D58").Value = 44364
$ws.Range("J58").Value = 80

# Row 59 (source row 51)
$ws.Range("D59").Value = 44355
$ws.Range("J59").Value = 240

# Row 60 (source row 100)
$ws.Range("D60").Value = 44418
$ws.Range("J60").Value = 250

# Row 61 (source row 8)
$ws.Range("D61").Value = 44278
$ws.Range("J61").Value = 250
$ws.Range("M61").Value = 14400
$ws.Range("P61").Value = 1440

# Row 62 (source row 63)
$ws.Range("D62").Value = 44245
$ws.Range("J62").Value = 120

# Row 63 (source row 43)
$ws.Range("D63").Value = 44230
$ws.Range("J63").Value = 80

# Row 64 (source row 5)
$ws.Range("D64").Value = 44382
$ws.Range("J64").Value = 70

# Row 65 (source row 58)
$ws.Range("D65").Value = 44284
$ws.Range("J65").Value = 40

# Row 66 (source row 47)
$ws.Range("D66").Value = 44449
$ws.Range("J66").Value = 220
$ws.Range("K66").Value = 18000
$ws.Range("L66").Value = 18000
$ws.Range("M66").Value = 18000
$ws.Range("P66").Value = 1800

# Row 67 (source row 82)
$ws.Range("D67").Value = 44336
$ws.Range("J67").Value = 120
$ws.Range("K67").Value = 14500
$ws.Range("M67").Value = 14750
$ws.Range("P67").Value = 1475

# Row 68 (source row 44)
$ws.Range("D68").Value = 44271
$ws.Range("J68").Value = 250
$ws.Range("M68").Value = 14400
$ws.Range("P68").Value = 1440

# Row 69 (source row 115)
$ws.Range("D69").Value = 44308
$ws.Range("J69").Value = 100
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = 15000
$ws.Range("P69").Value = 1500

# Row 70 (source row 33)
$ws.Range("D70").Value = 44330

# Row 71 (source row 117)
$ws.Range("D71").Value = 44463
$ws.Range("J71").Value = 240
$ws.Range("K71").Value = 17000
$ws.Range("L71").Value = 17000
$ws.Range("M71").Value = 17000
$ws.Range("P71").Value = 1700

# Row 72 (source row 48)
$ws.Range("D72").Value = 44322
$ws.Range("J72").Value = 80

# Row 73 (source row 17)
$ws.Range("D73").Value = 44428
$ws.Range("J73").Value = 250

# Row 74 (source row 40)
$ws.Range("D74").Value = 44224
$ws.Range("J74").Value = 200
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = 14500
$ws.Range("P74").Value = 1450

# Row 75 (source row 116)
$ws.Range("D75").Value = 44264
$ws.Range("J75").Value = 250
$ws.Range("K75").Value = 14500
$ws.Range("M75").Value = 14700
$ws.Range("P75").Value = 1470

# Row 76 (source row 11)
$ws.Range("D76").Value = 44315
$ws.Range("J76").Value = 120
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = 14500
$ws.Range("P76").Value = 1450

# Row 77 (source row 50)
$ws.Range("D77").Value = 44293
$ws.Range("J77").Value = 40
$ws.Range("K77").Value = 14000
$ws.Range("L77").Value = 14000
$ws.Range("M77").Value = 14000
$ws.Range("P77").Value = 1400

# Row 78 (source row 23)
$ws.Range("D78").Value = 44421
$ws.Range("J78").Value = 240
$ws.Range("K78").Value = 16500
$ws.Range("L78").Value = 17000
$ws.Range("M78").Value = 16750
$ws.Range("P78").Value = 1675

# Row 79 (source row 10)
$ws.Range("D79").Value = 44294
$ws.Range("J79").Value = 150
$ws.Range("K79").Value = 14000
$ws.Range("M79").Value = 14533
$ws.Range("P79").Value = 1453

# Row 80 (source row 36)
$ws.Range("D80").Value = 44232
$ws.Range("J80").Value = 300
$ws.Range("K80").Value = 14000
$ws.Range("M80").Value = 14500
$ws.Range("P80").Value = 1450

# Row 81 (source row 126)
$ws.Range("D81").Value = 44316
$ws.Range("J81").Value = 250
$ws.Range("K81").Value = 15000
$ws.Range("L81").Value = 15000
$ws.Range("M81").Value = 15000
$ws.Range("P81").Value = 1500

# Row 82 (source row 42)
$ws.Range("D82").Value = 44169
$ws.Range("J82").Value = 250
$ws.Range("K82").Value = 12000
$ws.Range("L82").Value = 12000
$ws.Range("M82").Value = 12000
$ws.Range("P82").Value = 1200

# Row 83 (source row 142)
$ws.Range("D83").Value = 44410
$ws.Range("J83").Value = 120

# Row 84 (source row 66)
$ws.Range("D84").Value = 44467
$ws.Range("K84").Value = 17000
$ws.Range("L84").Value = 17000
$ws.Range("M84").Value = 17000
$ws.Range("P84").Value = 1700

# Row 85 (source row 104)
$ws.Range("D85").Value = 44427
$ws.Range("J85").Value = 100
$ws.Range("K85").Value = 18000
$ws.Range("L85").Value = 18000
$ws.Range("M85").Value = 18000
$ws.Range("P85").Value = 1800

# Row 86 (source row 46)
$ws.Range("D86").Value = 44319
$ws.Range("J86").Value = 40

# Row 87 (source row 130)
$ws.Range("D87").Value = 44273
$ws.Range("J87").Value = 60
$ws.Range("K87").Value = 15000
$ws.Range("L87").Value = 15000
$ws.Range("M87").Value = 15000
$ws.Range("P87").Value = 1500

# Row 88 (source row 143)
$ws.Range("D88").Value = 44411

# Row 89 (source row 21)
$ws.Range("D89").Value = 44236
$ws.Range("J89").Value = 250
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = 15000
$ws.Range("P89").Value = 1500

# Row 90 (source row 137)
$ws.Range("D90").Value = 44215

# Row 91 (source row 14)
$ws.Range("D91").Value = 44403
$ws.Range("J91").Value = 40
$ws.Range("K91").Value = 15000
$ws.Range("L91").Value = 15000
$ws.Range("M91").Value = 15000
$ws.Range("P91").Value = 1500

# Row 92 (source row 140)
$ws.Range("D92").Value = 44203
$ws.Range("J92").Value = 100
$ws.Range("K92").Value = 12500
$ws.Range("L92").Value = 12500
$ws.Range("M92").Value = 12500
$ws.Range("P92").Value = 1250

# Row 93 (source row 3)
$ws.Range("D93").Value = 44413
$ws.Range("J93").Value = 150
$ws.Range("K93").Value = 15000
$ws.Range("M93").Value = 15000
$ws.Range("P93").Value = 1500

# Row 94 (source row 54)
$ws.Range("D94").Value = 44469
$ws.Range("H94").Value = "Chilote"
$ws.Range("J94").Value = 120
$ws.Range("K94").Value = 17000
$ws.Range("L94").Value = 17000
$ws.Range("M94").Value = 17000
$ws.Range("P94").Value = 1700

# Row 95 (source row 123)
$ws.Range("D95").Value = 44168
$ws.Range("J95").Value = 100
$ws.Range("K95").Value = 12000
$ws.Range("L95").Value = 12000
$ws.Range("M95").Value = 12000
$ws.Range("P95").Value = 1200

# Row 96 (source row 69)
$ws.Range("D96").Value = 44434
$ws.Range("J96").Value = 120
$ws.Range("K96").Value = 18000
$ws.Range("L96").Value = 18000
$ws.Range("M96").Value = 18000
$ws.Range("P96").Value = 1800

# Row 97 (source row 148)
$ws.Range("D97").Value = 44292
$ws.Range("J97").Value = 250
$ws.Range("L97").Value = 14000
$ws.Range("M97").Value = 14000
$ws.Range("P97").Value = 1400

# Row 98 (source row 78)
$ws.Range("D98").Value = 44161
$ws.Range("K98").Value = 12000
$ws.Range("L98").Value = 12000
$ws.Range("M98").Value = 12000
$ws.Range("P98").Value = 1200

# Row 99 (source row 134)
$ws.Range("D99").Value = 44351
$ws.Range("J99").Value = 220

# Row 100 (source row 109)
$ws.Range("D100").Value = 44280
$ws.Range("J100").Value = 100
$ws.Range("K100").Value = 14000
$ws.Range("M100").Value = 14500
$ws.Range("P100").Value = 1450

# Row 101 (source row 84)
$ws.Range("D101").Value = 44274
$ws.Range("J101").Value = 250
$ws.Range("K101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("M101").Value = 15000
$ws.Range("P101").Value = 1500

# Row 102 (source row 61)
$ws.Range("D102").Value = 44229
$ws.Range("J102").Value = 300
$ws.Range("K102").Value = 14000
$ws.Range("M102").Value = 14500
$ws.Range("P102").Value = 1450

# Row 103 (source row 96)
$ws.Range("D103").Value = 44344
$ws.Range("J103").Value = 240
$ws.Range("K103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 15000
$ws.Range("P103").Value = 1500

# Row 104 (source row 79)
$ws.Range("D104").Value = 44407
$ws.Range("J104").Value = 240
$ws.Range("K104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("M104").Value = 15000
$ws.Range("P104").Value = 1500

# Row 105 (source row 149)
$ws.Range("D105").Value = 44358
$ws.Range("J105").Value = 240
$ws.Range("K105").Value = 15000
$ws.Range("L105").Value = 15000
$ws.Range("M105").Value = 15000
$ws.Range("P105").Value = 1500

# Row 106 (source row 7)
$ws.Range("D106").Value = 44335
$ws.Range("J106").Value = 40
$ws.Range("K106").Value = 15000
$ws.Range("L106").Value = 15000
$ws.Range("M106").Value = 15000
$ws.Range("P106").Value = 1500

# Row 107 (source row 103)
$ws.Range("D107").Value = 44165
$ws.Range("J107").Value = 80
$ws.Range("K107").Value = 12000
$ws.Range("L107").Value = 12000
$ws.Range("M107").Value = 12000
$ws.Range("P107").Value = 1200

# Row 108 (source row 37)
$ws.Range("D108").Value = 44231
$ws.Range("J108").Value = 80
$ws.Range("K108").Value = 14000
$ws.Range("L108").Value = 15000
$ws.Range("M108").Value = 14500
$ws.Range("P108").Value = 1450

# Row 109 (source row 56)
$ws.Range("D109").Value = 44320
$ws.Range("J109").Value = 240
$ws.Range("K109").Value = 15000
$ws.Range("M109").Value = 15000
$ws.Range("P109").Value = 1500

# Row 110 (source row 121)
$ws.Range("D110").Value = 44252
$ws.Range("J110").Value = 80
$ws.Range("K110").Value = 15000
$ws.Range("M110").Value = 15000
$ws.Range("P110").Value = 1500

# Row 111 (source row 6)
$ws.Range("D111").Value = 44204
$ws.Range("J111").Value = 240
$ws.Range("K111").Value = 13000
$ws.Range("L111").Value = 13000
$ws.Range("M111").Value = 13000
$ws.Range("P111").Value = 1300

# Row 112 (source row 18)
$ws.Range("D112").Value = 44362
$ws.Range("J112").Value = 240

# Row 113 (source row 71)
$ws.Range("D113").Value = 44300
$ws.Range("J113").Value = 20

# Row 114 (source row 15)
$ws.Range("D114").Value = 44186
$ws.Range("J114").Value = 80
$ws.Range("K114").Value = 13000
$ws.Range("L114").Value = 13000
$ws.Range("M114").Value = 13000
$ws.Range("P114").Value = 1300

# Row 115 (source row 132)
$ws.Range("D115").Value = 44372
$ws.Range("J115").Value = 200

# Row 116 (source row 99)
$ws.Range("D116").Value = 44350
$ws.Range("J116").Value = 70
$ws.Range("K116").Value = 15000
$ws.Range("M116").Value = 15000
$ws.Range("P116").Value = 1500

# Row 117 (source row 91)
$ws.Range("D117").Value = 44448
$ws.Range("J117").Value = 120
$ws.Range("K117").Value = 18000
$ws.Range("L117").Value = 18000
$ws.Range("M117").Value = 18000
$ws.Range("P117").Value = 1800

# Row 118 (source row 120)
$ws.Range("D118").Value = 44243
$ws.Range("J118").Value = 250
$ws.Range("K118").Value = 15000
$ws.Range("L118").Value = 15000
$ws.Range("M118").Value = 15000
$ws.Range("P118").Value = 1500

# Row 119 (source row 57)
$ws.Range("D119").Value = 44385

# Row 120 (source row 73)
$ws.Range("D120").Value = 44452
$ws.Range("J120").Value = 80
$ws.Range("K120").Value = 18000
$ws.Range("L120").Value = 18000
$ws.Range("M120").Value = 18000
$ws.Range("P120").Value = 1800

# Row 121 (source row 20)
$ws.Range("D121").Value = 44202
$ws.Range("J121").Value = 50
$ws.Range("K121").Value = 13000
$ws.Range("L121").Value = 13000
$ws.Range("M121").Value = 13000
$ws.Range("P121").Value = 1300

# Row 122 (source row 39)
$ws.Range("D122").Value = 44435
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = 18000
$ws.Range("P122").Value = 1800

# Row 123 (source row 31)
$ws.Range("D123").Value = 44242
$ws.Range("J123").Value = 80
$ws.Range("K123").Value = 15000
$ws.Range("L123").Value = 15000
$ws.Range("M123").Value = 15000
$ws.Range("P123").Value = 1500

# Row 124 (source row 35)
$ws.Range("D124").Value = 44377
$ws.Range("J124").Value = 30

# Row 125 (source row 80)
$ws.Range("D125").Value = 44307
$ws.Range("J125").Value = 20

# Row 126 (source row 32)
$ws.Range("D126").Value = 44159
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = 12000
$ws.Range("P126").Value = 1200

# Row 127 (source row 122)
$ws.Range("D127").Value = 44166
$ws.Range("K127").Value = 12000
$ws.Range("L127").Value = 12000
$ws.Range("M127").Value = 12000
$ws.Range("P127").Value = 1200

# Row 128 (source row 64)
$ws.Range("D128").Value = 44334
$ws.Range("J128").Value = 240

# Row 129 (source row 114)
$ws.Range("D129").Value = 44386
$ws.Range("J129").Value = 240
$ws.Range("K129").Value = 15000
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = 15000
$ws.Range("P129").Value = 1500

# Row 130 (source row 26)
$ws.Range("D130").Value = 44384
$ws.Range("J130").Value = 15

# Row 131 (source row 128)
$ws.Range("D131").Value = 44363
$ws.Range("J131").Value = 20
$ws.Range("K131").Value = 15000
$ws.Range("L131").Value = 15000
$ws.Range("M131").Value = 15000
$ws.Range("P131").Value = 1500

# Row 132 (source row 136)
$ws.Range("D132").Value = 44306
$ws.Range("J132").Value = 250

# Row 133 (source row 124)
$ws.Range("D133").Value = 44369
$ws.Range("J133").Value = 250
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = 15000
$ws.Range("P133").Value = 1500

# Row 134 (source row 105)
$ws.Range("D134").Value = 44172
$ws.Range("J134").Value = 100
$ws.Range("K134").Value = 12000
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = 12000
$ws.Range("P134").Value = 1200

# Row 135 (source row 119)
$ws.Range("D135").Value = 44301
$ws.Range("J135").Value = 80

# Row 136 (source row 139)
$ws.Range("D136").Value = 44357
$ws.Range("J136").Value = 90

# Row 137 (source row 76)
$ws.Range("D137").Value = 44223
$ws.Range("J137").Value = 60

# Row 138 (source row 45)
$ws.Range("D138").Value = 44298
$ws.Range("J138").Value = 40
$ws.Range("K138").Value = 13000
$ws.Range("L138").Value = 13000
$ws.Range("M138").Value = 13000
$ws.Range("P138").Value = 1300

# Row 139 (source row 127)
$ws.Range("D139").Value = 44397
$ws.Range("J139").Value = 250

# Row 140 (source row 102)
$ws.Range("D140").Value = 44414
$ws.Range("J140").Value = 280
$ws.Range("K140").Value = 15000
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = 15000
$ws.Range("P140").Value = 1500

# Row 141 (source row 49)
$ws.Range("D141").Value = 44333
$ws.Range("J141").Value = 40
$ws.Range("K141").Value = 15000
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 15000
$ws.Range("P141").Value = 1500

# Row 142 (source row 86)
$ws.Range("D142").Value = 44392
$ws.Range("J142").Value = 240

# Row 143 (source row 151)
$ws.Range("D143").Value = 44425
$ws.Range("K143").Value = 18000
$ws.Range("L143").Value = 18000
$ws.Range("M143").Value = 18000
$ws.Range("P143").Value = 1800

# Row 144 (source row 113)
$ws.Range("D144").Value = 44390
$ws.Range("J144").Value = 240

# Row 145 (source row 106)
$ws.Range("D145").Value = 44466
$ws.Range("J145").Value = 180
$ws.Range("K145").Value = 17000
$ws.Range("L145").Value = 17000
$ws.Range("M145").Value = 17000
$ws.Range("P145").Value = 1700

# Row 146 (source row 65)
$ws.Range("D146").Value = 44250
$ws.Range("J146").Value = 300

# Row 147 (source row 131)
$ws.Range("D147").Value = 44438
$ws.Range("J147").Value = 100
$ws.Range("K147").Value = 18000
$ws.Range("L147").Value = 18000
$ws.Range("M147").Value = 18000
$ws.Range("P147").Value = 1800

# Row 148 (source row 108)
$ws.Range("D148").Value = 44201
$ws.Range("J148").Value = 240
$ws.Range("K148").Value = 12500
$ws.Range("L148").Value = 13000
$ws.Range("M148").Value = 12750
$ws.Range("P148").Value = 1275

# Row 149 (source row 81)
$ws.Range("D149").Value = 44193
$ws.Range("J149").Value = 180
$ws.Range("K149").Value = 13000
$ws.Range("L149").Value = 13500
$ws.Range("M149").Value = 13250
$ws.Range("P149").Value = 1325

# Row 150 (source row 112)
$ws.Range("D150").Value = 44312
$ws.Range("J150").Value = 25
$ws.Range("K150").Value = 15000
$ws.Range("M150").Value = 15000
$ws.Range("P150").Value = 1500

# Row 151 (source row 16)
$ws.Range("D151").Value = 44326
$ws.Range("J151").Value = 40
$ws.Range("K151").Value = 15000
$ws.Range("L151").Value = 15000
$ws.Range("M151").Value = 15000
$ws.Range("P151").Value = 1500

# Row 152 (source row 9)
$ws.Range("D152").Value = 44432
$ws.Range("J152").Value = 240
$ws.Range("K152").Value = 18000
$ws.Range("L152").Value = 18000
$ws.Range("M152").Value = 18000
$ws.Range("P152").Value = 1800
